$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("展览")
$ws.Range("F2").Value = 37984
$ws.Range("F5").Value = 800
$ws.Range("F7").Value = 378
$ws.Range("F8").Value = 471
$ws.Range("F9").Value = 869
$ws.Range("F10").Value = 110
$ws.Range("F11").Value = 755
$ws.Range("F12").Value = 591
$ws.Range("F13").Value = 85
$ws.Range("F14").Value = 40
$ws.Range("F15").Value = 43
$ws.Range("F16").Value = 691
$ws.Range("F17").Value = 190
$ws.Range("F18").Value = 493
$ws.Range("F19").Value = 453
$ws.Range("F20").Value = 1194
$ws.Range("F21").Value = 99
$ws.Range("F22").Value = 884
$ws.Range("F23").Value = 2601
$ws.Range("F24").Value = 1086
$ws.Range("F27").Value = 1182
$ws.Range("F29").Value = 840
$ws.Range("F31").Value = 1187

$ws = $wb.Worksheets.Item("演出")
$ws.Range("F3").Value = 450

$ws = $wb.Worksheets.Item("本地生活")
$ws.Range("F2").Value = 674

$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F2").Value = 674
$ws.Range("F3").Value = 37984
$ws.Range("F6").Value = 800
$ws.Range("F9").Value = 378
$ws.Range("F10").Value = 471
$ws.Range("F11").Value = 450
$ws.Range("F12").Value = 450
$ws.Range("F16").Value = 869
$ws.Range("F17").Value = 110
$ws.Range("F18").Value = 755
$ws.Range("F19").Value = 591
$ws.Range("F20").Value = 85
$ws.Range("F22").Value = 40
$ws.Range("F26").Value = 43
$ws.Range("F28").Value = 691
$ws.Range("F29").Value = 190
$ws.Range("F30").Value = 493
$ws.Range("F31").Value = 453
$ws.Range("F32").Value = 1194
$ws.Range("F33").Value = 99
$ws.Range("F34").Value = 884
$ws.Range("F35").Value = 2601
$ws.Range("F36").Value = 1086
$ws.Range("F39").Value = 1182
$ws.Range("F42").Value = 840
$ws.Range("F44").Value = 1187
